# Applies the "update final UI of notbook" edit:
#  1. Append city/country location info to each numbered affiliation line.
#  2. Re-shuffle which initials list sits under which "Author Contributions"
#     heading (Cohort PI / Cohort co-investigator / Performed the analysis /
#     Read, edited and approved the paper / Contributed data or analysis
#     tools / Analyzed the data all swap content with one another).

$d = $word.ActiveDocument

function Replace-InDoc($old, $new) {
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

function Replace-InParagraph($index, $old, $new) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $null = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- 1. Affiliation list: append ", City, Country" style suffixes ---------

Replace-InDoc "Imaging Genetics Center, Stevens Neuroimaging & Informatics Institute, USC" "Imaging Genetics Center, Stevens Neuroimaging & Informatics Institute, USC, Los Angeles, CA, United States"
Replace-InDoc "Amsterdam UMC" "Amsterdam UMC, Amsterdam, Netherlands"
Replace-InDoc "Stanford University" "Stanford University, Palo Alto, CA, United States"
Replace-InDoc "IRCCS Fondazione Santa Lucia" "IRCCS Fondazione Santa Lucia, Rome, Italy"
Replace-InDoc "Departamento de Neurologia, University of Campinas (UNICAMP)" "Departamento de Neurologia, University of Campinas (UNICAMP), Campinas, Brazil"
Replace-InDoc "Support Centre for Advanced Neuroimaging (SCAN), Inselspital" "Support Centre for Advanced Neuroimaging (SCAN), Inselspital, Bern, Switzerland"
Replace-InDoc "Department of Neurology, Inselspital, University of Bern" "Department of Neurology, Inselspital, University of Bern, Bern, Switzerland"
Replace-InDoc "Institute of Diagnostic and Interventional Neuroradiology, Inselspital " "Institute of Diagnostic and Interventional Neuroradiology, Inselspital , Bern, Switzerland"
Replace-InDoc "Division of Neuroscience and Experimental Psychology, University of Manchester" "Division of Neuroscience and Experimental Psychology, University of Manchester, Manchester, United Kingdom"
Replace-InDoc "Division of Neuroradiology, University of Virginia" "Division of Neuroradiology, University of Virginia, Charlottesville, VA, United States"
Replace-InDoc "Illinois Institute of Technology" "Illinois Institute of Technology, Chicago, United States"
Replace-InDoc "Department of Medicine, University of Otago" "Department of Medicine, University of Otago, Christchurch, New Zealand"
Replace-InDoc "Fondazione IRCCS" "Fondazione IRCCS, Milan, Italy"
Replace-InDoc "University of Pennsylvania, Penn Frontotemporal Degeneration Center" "University of Pennsylvania, Penn Frontotemporal Degeneration Center, Philadelphia, PA, United States"

# --- 2. Author Contributions: swap which list goes under which heading ----
# Locate the nine contribution paragraphs (after the "Author Contributions"
# heading) once, by their current (pre-edit) title text, so the swap is
# robust even if the heading/paragraph numbers were to shift slightly.

$contribStart = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "Author Contributions (in alphabetical order)") {
        $contribStart = $i
        break
    }
}

$idxPerformed   = $contribStart + 1   # "Performed the analysis"
$idxWrote       = $contribStart + 2   # "Wrote the paper" (unchanged)
$idxContributed = $contribStart + 3   # "Contributed data or analysis tools"
$idxAnalyzed    = $contribStart + 4   # "Analyzed the data"
$idxConceived   = $contribStart + 5   # "Conceived and designed the analysis" (unchanged)
$idxCollected   = $contribStart + 6   # "Collected the data" (unchanged)
$idxCohortPI    = $contribStart + 7   # "Cohort PI"
$idxReadEdited  = $contribStart + 8   # "Read, edited and approved the paper"
$idxCohortCoInv = $contribStart + 9   # "Cohort co-investigator"

# Content that will move into each paragraph slot (title + initials list),
# captured as literal strings so the six-way cycle below applies cleanly.

$performedTitle   = "Performed the analysis"
$performedList    = "J.B., M.L., Y.D.vdW."

$contributedTitle = "Contributed data or analysis tools"
$contributedList  = "J.B., M.L., Y.D.vdW."

$analyzedTitle    = "Analyzed the data"
$analyzedList     = "C.R., C.Y., F.C., F.P., G.S., J.B., J.D., K.L.P., K.Z., M.L., M.R., T.M., T.P., Y.D.vdW."

$cohortPITitle    = "Cohort PI"
$cohortPIList     = "Y.D.vdW."

$readEditedTitle  = "Read, edited and approved the paper"
$readEditedList   = "B.G., C.M., C.R., C.R., C.Y., F.C., F.P., G.S., I.D., J.D., K.L.P., K.Z., L.P., M.R., N.J., O.A.vdH., P.T., R.M.DB., R.W., S.a-B., T.M., T.P., Y.D.vdW."

$cohortCoInvTitle = "Cohort co-investigator"
$cohortCoInvList  = "J.B., M.L., N.J., P.T."

# Slot 1 (was "Performed the analysis") <- "Cohort PI" content
Replace-InParagraph $idxPerformed   $performedTitle   $cohortPITitle
Replace-InParagraph $idxPerformed   $performedList    $cohortPIList

# Slot 3 (was "Contributed data or analysis tools") <- "Cohort co-investigator" content
Replace-InParagraph $idxContributed $contributedTitle $cohortCoInvTitle
Replace-InParagraph $idxContributed $contributedList  $cohortCoInvList

# Slot 4 (was "Analyzed the data") <- "Performed the analysis" content
Replace-InParagraph $idxAnalyzed    $analyzedTitle    $performedTitle
Replace-InParagraph $idxAnalyzed    $analyzedList     $performedList

# Slot 7 (was "Cohort PI") <- "Read, edited and approved the paper" content
Replace-InParagraph $idxCohortPI    $cohortPITitle    $readEditedTitle
Replace-InParagraph $idxCohortPI    $cohortPIList     $readEditedList

# Slot 8 (was "Read, edited and approved the paper") <- "Contributed data or analysis tools" content
Replace-InParagraph $idxReadEdited  $readEditedTitle  $contributedTitle
Replace-InParagraph $idxReadEdited  $readEditedList   $contributedList

# Slot 9 (was "Cohort co-investigator") <- "Analyzed the data" content
Replace-InParagraph $idxCohortCoInv $cohortCoInvTitle $analyzedTitle
Replace-InParagraph $idxCohortCoInv $cohortCoInvList  $analyzedList

Write-Output "done"
